$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129 — this shifts the existing rows
# 129..218 down to 130..219 (and the dimension grows to A1:R219).
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new data point.
$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44574
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = 100112037
$ws.Range("G129").Value = "Cebollín"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 70
$ws.Range("K129").Value = 6500
$ws.Range("L129").Value = 6500
$ws.Range("M129").Value = 6500
$ws.Range("N129").Value = "`$/paquete 36 unidades"
$ws.Range("O129").Value = "Región Metropolitana"
$ws.Range("P129").Value = 181
$ws.Range("Q129").Value = 36
$ws.Range("R129").Value = "Hortaliza"
